$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New monthly data rows (22-31), mirroring the existing Month/Pageviews columns.
$dates = @(42705,42736,42767,42795,42826,42856,42887,42917,42948,42979)
$views = @(10051,8004,6080,6444,5994,5598,1800,2497,2410,2210)

for ($i = 0; $i -lt $dates.Length; $i++) {
    $row = 22 + $i
    $aCell = $ws.Cells.Item($row, 1)
    $aCell.Value = $dates[$i]
    $aCell.NumberFormat = "mmm-yy"

    $bCell = $ws.Cells.Item($row, 2)
    $bCell.Value = $views[$i]
}

# B22 carries an explicit "#,##0" number format (new style), unlike the other B cells.
$ws.Cells.Item(22, 2).NumberFormat = "#,##0"

# Row 32: A32 only, formatted like the date column but left empty.
$a32 = $ws.Cells.Item(32, 1)
$a32.NumberFormat = "mmm-yy"

[void]$ws.Range("E30:E31").Select()
